$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild data rows (2-13) with the updated TPM values: refreshed numbers for the
# existing FAPs->* rows (2-7) plus six new MuSCs->* rows (8-13) for the same
# Fgl1/Lag3 ligand-receptor pair.
$ws.Range("A2:T13").Clear()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgl1"
$ws.Range("C2").Value = "Lag3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.085463
$ws.Range("H2").Value = 0.256389
$ws.Range("I2").Value = 0.6706908097814145
$ws.Range("J2").Value = 0.6706908097814145
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.806827333333334
$ws.Range("N2").Value = 23.420482
$ws.Range("O2").Value = 0.3564356619476692
$ws.Range("P2").Value = 0.3564356619476692
$ws.Range("Q2").Value = 0.6671948843886667
$ws.Range("R2").Value = 6.004753959497999
$ws.Range("S2").Value = 0.2390581227466568
$ws.Range("T2").Value = 0.2390581227466567

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgl1"
$ws.Range("C3").Value = "Lag3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.085463
$ws.Range("H3").Value = 0.256389
$ws.Range("I3").Value = 0.6706908097814145
$ws.Range("J3").Value = 0.6706908097814145
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.661742
$ws.Range("N3").Value = 13.985226
$ws.Range("O3").Value = 0.2128407641993771
$ws.Range("P3").Value = 0.212840764199377
$ws.Range("Q3").Value = 0.398406456546
$ws.Range("R3").Value = 3.585658108914
$ws.Range("S3").Value = 0.1427503444953753
$ws.Range("T3").Value = 0.1427503444953753

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgl1"
$ws.Range("C4").Value = "Lag3"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.085463
$ws.Range("H4").Value = 0.256389
$ws.Range("I4").Value = 0.6706908097814145
$ws.Range("J4").Value = 0.6706908097814145
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.225227333333333
$ws.Range("N4").Value = 6.675682
$ws.Range("O4").Value = 0.1015970180554841
$ws.Range("P4").Value = 0.1015970180554841
$ws.Range("Q4").Value = 0.1901746035886666
$ws.Range("R4").Value = 1.711571432298
$ws.Range("S4").Value = 0.06814018631100963
$ws.Range("T4").Value = 0.06814018631100963

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgl1"
$ws.Range("C5").Value = "Lag3"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.085463
$ws.Range("H5").Value = 0.256389
$ws.Range("I5").Value = 0.6706908097814145
$ws.Range("J5").Value = 0.6706908097814145
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.309581
$ws.Range("N5").Value = 6.928743
$ws.Range("O5").Value = 0.1054483463521494
$ws.Range("P5").Value = 0.1054483463521493
$ws.Range("Q5").Value = 0.197383721003
$ws.Range("R5").Value = 1.776453489027
$ws.Range("S5").Value = 0.07072323680503412
$ws.Range("T5").Value = 0.0707232368050341

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgl1"
$ws.Range("C6").Value = "Lag3"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.085463
$ws.Range("H6").Value = 0.256389
$ws.Range("I6").Value = 0.6706908097814145
$ws.Range("J6").Value = 0.6706908097814145
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.936864333333334
$ws.Range("N6").Value = 8.810593
$ws.Range("O6").Value = 0.134088168984161
$ws.Range("P6").Value = 0.134088168984161
$ws.Range("Q6").Value = 0.2509932365196667
$ws.Range("R6").Value = 2.258939128677
$ws.Range("S6").Value = 0.0899317026380941
$ws.Range("T6").Value = 0.08993170263809408

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgl1"
$ws.Range("C7").Value = "Lag3"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.085463
$ws.Range("H7").Value = 0.256389
$ws.Range("I7").Value = 0.6706908097814145
$ws.Range("J7").Value = 0.6706908097814145
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.962244666666667
$ws.Range("N7").Value = 5.886734000000001
$ws.Range("O7").Value = 0.0895900404611592
$ws.Range("P7").Value = 0.0895900404611592
$ws.Range("Q7").Value = 0.1676993159473333
$ws.Range("R7").Value = 1.509293843526
$ws.Range("S7").Value = 0.06008721678524456
$ws.Range("T7").Value = 0.06008721678524456

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Fgl1"
$ws.Range("C8").Value = "Lag3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.04196233333333333
$ws.Range("H8").Value = 0.125887
$ws.Range("I8").Value = 0.3293091902185855
$ws.Range("J8").Value = 0.3293091902185856
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.806827333333334
$ws.Range("N8").Value = 23.420482
$ws.Range("O8").Value = 0.3564356619476692
$ws.Range("P8").Value = 0.3564356619476692
$ws.Range("Q8").Value = 0.3275926908371111
$ws.Range("R8").Value = 2.948334217534
$ws.Range("S8").Value = 0.1173775392010124
$ws.Range("T8").Value = 0.1173775392010124

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Fgl1"
$ws.Range("C9").Value = "Lag3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.04196233333333333
$ws.Range("H9").Value = 0.125887
$ws.Range("I9").Value = 0.3293091902185855
$ws.Range("J9").Value = 0.3293091902185856
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.661742
$ws.Range("N9").Value = 13.985226
$ws.Range("O9").Value = 0.2128407641993771
$ws.Range("P9").Value = 0.212840764199377
$ws.Range("Q9").Value = 0.195617571718
$ws.Range("R9").Value = 1.760558145462
$ws.Range("S9").Value = 0.07009041970400176
$ws.Range("T9").Value = 0.07009041970400177

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgl1"
$ws.Range("C10").Value = "Lag3"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.04196233333333333
$ws.Range("H10").Value = 0.125887
$ws.Range("I10").Value = 0.3293091902185855
$ws.Range("J10").Value = 0.3293091902185856
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.225227333333333
$ws.Range("N10").Value = 6.675682
$ws.Range("O10").Value = 0.1015970180554841
$ws.Range("P10").Value = 0.1015970180554841
$ws.Range("Q10").Value = 0.09337573110377777
$ws.Range("R10").Value = 0.840381579934
$ws.Range("S10").Value = 0.03345683174447448
$ws.Range("T10").Value = 0.03345683174447449

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Fgl1"
$ws.Range("C11").Value = "Lag3"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.04196233333333333
$ws.Range("H11").Value = 0.125887
$ws.Range("I11").Value = 0.3293091902185855
$ws.Range("J11").Value = 0.3293091902185856
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.309581
$ws.Range("N11").Value = 6.928743
$ws.Range("O11").Value = 0.1054483463521494
$ws.Range("P11").Value = 0.1054483463521493
$ws.Range("Q11").Value = 0.09691540778233333
$ws.Range("R11").Value = 0.872238670041
$ws.Range("S11").Value = 0.03472510954711524
$ws.Range("T11").Value = 0.03472510954711524

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Fgl1"
$ws.Range("C12").Value = "Lag3"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.04196233333333333
$ws.Range("H12").Value = 0.125887
$ws.Range("I12").Value = 0.3293091902185855
$ws.Range("J12").Value = 0.3293091902185856
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.936864333333334
$ws.Range("N12").Value = 8.810593
$ws.Range("O12").Value = 0.134088168984161
$ws.Range("P12").Value = 0.134088168984161
$ws.Range("Q12").Value = 0.1232376801101111
$ws.Range("R12").Value = 1.109139120991
$ws.Range("S12").Value = 0.04415646634606692
$ws.Range("T12").Value = 0.04415646634606692

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Fgl1"
$ws.Range("C13").Value = "Lag3"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.04196233333333333
$ws.Range("H13").Value = 0.125887
$ws.Range("I13").Value = 0.3293091902185855
$ws.Range("J13").Value = 0.3293091902185856
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.962244666666667
$ws.Range("N13").Value = 5.886734000000001
$ws.Range("O13").Value = 0.0895900404611592
$ws.Range("P13").Value = 0.0895900404611592
$ws.Range("Q13").Value = 0.08234036478422223
$ws.Range("R13").Value = 0.7410632830580001
$ws.Range("S13").Value = 0.02950282367591465
$ws.Range("T13").Value = 0.02950282367591466
